$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark that sits right after <trinh_ky>.
# ------------------------------------------------------------------
$oldGoBack = $d.Bookmarks.Item("_GoBack")
$oldGoBack.Delete()

# ------------------------------------------------------------------
# 2) Locate "ngày 13 tháng 6 năm 2019" inside the "Căn cứ Luật Quản lý
#    thuế ..." sentence and replace it with the "<luat_qlt_ngay>"
#    placeholder.
# ------------------------------------------------------------------
$matchRange = $d.Content
$matchRange.Find.Execute("ngày 13 tháng 6 năm 2019", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)
$matchStart = $matchRange.Start
$matchEnd = $matchRange.End

$matchRange.Text = "<luat_qlt_ngay>"
$placeholderEnd = $matchRange.End

# Split the run right before the placeholder text so the preceding
# "n cứ Luật Quản lý thuế " text keeps its own run (toggling a format
# on/off forces the run boundary without leaving any stray formatting,
# and - unlike inserting/removing a bookmark at that spot - correctly
# recomputes xml:space on each side).
$placeholderRange = $d.Range($matchStart, $placeholderEnd)
$placeholderRange.Font.Bold = 1
$placeholderRange.Font.Bold = 0

# ------------------------------------------------------------------
# 3) Re-create the "_GoBack" bookmark right after the new placeholder.
# ------------------------------------------------------------------
$bookmarkRange = $d.Range($placeholderEnd, $placeholderEnd)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# ------------------------------------------------------------------
# 4) Split the trailing " và các văn bản hướng dẫn thi hành;" text so
#    the single space and the rest of the sentence become separate
#    runs, matching the target structure.
# ------------------------------------------------------------------
$tailStart = $placeholderEnd + 1
$tailEnd = $d.Content.End
$tailRange = $d.Range($tailStart, $tailEnd)
$tailRange.Find.Execute("và các văn bản hướng dẫn thi hành;", $true, $false, `
    $false, $false, $false, $true, 1, $false, "", 0)
$tailRange.Font.Bold = 1
$tailRange.Font.Bold = 0
